$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain
# exact-text (to match source formatting, avoid float coercion).
# Mark them as Text format BEFORE assigning the value.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D16", "D19", "D20", "D21", "D23", "D24", "D26", "D28", "D30", "D32", "D33", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume / coin data values.
$ws.Range("D2").Value = "58.218.04"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.473.11"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "520.78"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").Value = "134.31"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "2.483.13"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "0.0983"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "0.338"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "2.918.46"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "58.156.41"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "22.16"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").Value = "2.476.48"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.18"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "320.55"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  -5.56%  "
$ws.Range("D24").Value = "64.47"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").Value = "7.38"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "169.53"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "6.29"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "18.11"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "1.32"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "36.58"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Value = "0.801"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.16"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.45"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "273.62"
$ws.Range("E44").Value = "  -2.94%  "
$ws.Range("D45").Value = "0.596"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "124.23"
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("D47").Value = "0.0909"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").Value = "17.03"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "1.738.63"
$ws.Range("E51").Value = "  -1.07%  "
